# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets to the refreshed values captured at the later scrape.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2894
$ws1.Range("F3").Value = 21237
$ws1.Range("F4").Value = 103
$ws1.Range("F5").Value = 3079
$ws1.Range("F6").Value = 811
$ws1.Range("F8").Value = 519
$ws1.Range("F9").Value = 773
$ws1.Range("F10").Value = 283
$ws1.Range("F14").Value = 528
$ws1.Range("F16").Value = 273
$ws1.Range("F17").Value = 20
$ws1.Range("F18").Value = 425
$ws1.Range("F19").Value = 71
$ws1.Range("F20").Value = 27
$ws1.Range("F22").Value = 41

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 341

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 1644

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1644
$ws4.Range("F6").Value = 2894
$ws4.Range("F7").Value = 21237
$ws4.Range("F10").Value = 103
$ws4.Range("F12").Value = 341
$ws4.Range("F13").Value = 3079
$ws4.Range("F14").Value = 811
$ws4.Range("F18").Value = 519
$ws4.Range("F19").Value = 773
$ws4.Range("F20").Value = 283
$ws4.Range("F29").Value = 528
$ws4.Range("F33").Value = 273
$ws4.Range("F36").Value = 20
$ws4.Range("F37").Value = 425
$ws4.Range("F39").Value = 71
$ws4.Range("F40").Value = 27
$ws4.Range("F44").Value = 41
